$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update Residue Combination text
$ws.Range("A5").Value = "130, 426"

# Row 8: data from former row 15 merged into row 8
$ws.Range("A8").Value = "130, 786"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "786, 130"
$ws.Range("D8").Value = "959, 1219"

# Row 10
$ws.Range("A10").Value = "98, 786, 1082"

# Row 11
$ws.Range("A11").Value = "98, 1082, SF"

# Row 12
$ws.Range("A12").Value = "754, 786, 1082"

# Row 13 - value looks numeric ("458"); force text so it matches the
# original inline-string typing used throughout this column.
$cell = $ws.Range("A13")
$cell.NumberFormat = "@"
$cell.Value = "458"
$cell.Style = "Normal"

# Row 14 - same numeric-looking text issue ("786")
$cell = $ws.Range("A14")
$cell.NumberFormat = "@"
$cell.Value = "786"
$cell.Style = "Normal"

# Row 15 is removed entirely - its data was combined into row 8 above
$ws.Rows.Item(15).Delete()
